$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rounds sheet: add a new round row for CGC-W-06022025
# ---------------------------------------------------------------------------
$rounds = $wb.Worksheets.Item("Rounds")

$rounds.Cells.Item(6, 1).Value = "CGC-W-06022025"
$rounds.Cells.Item(6, 2).Value = "CGC-W"
$rounds.Cells.Item(6, 3).Value = 45810

# Match the date formatting used by the other rows in the Date column
$rounds.Cells.Item(5, 3).Copy()
$rounds.Cells.Item(6, 3).PasteSpecial(-4122)  # xlPasteFormats

$rounds.Range("C7").Select()

# ---------------------------------------------------------------------------
# 2) Existing round-detail sheets: insert a "Tee Fairway" column before the
#    "Fairway Hits" column (without disturbing the <cols> width definition,
#    which stays pinned to column C).
# ---------------------------------------------------------------------------

function Add-TeeFairwayColumn($ws, [bool]$styleHeader) {
    $oldC1 = $ws.Cells.Item(1, 3).Value2
    $oldD1 = $ws.Cells.Item(1, 4).Value2
    $oldE1 = $ws.Cells.Item(1, 5).Value2

    $ws.Cells.Item(1, 6).Value = $oldE1
    $ws.Cells.Item(1, 5).Value = $oldD1
    $ws.Cells.Item(1, 4).Value = $oldC1
    $ws.Cells.Item(1, 3).Value = "Tee Fairway"

    for ($r = 10; $r -ge 2; $r--) {
        $val = $ws.Cells.Item($r, 5).Value2
        if ($null -ne $val) {
            $ws.Cells.Item($r, 6).Value = $val
        }
        $ws.Cells.Item($r, 5).ClearContents()
    }

    if ($styleHeader) {
        $ws.Range("C1:F1").Font.Color = 0
    }
}

$ws2 = $wb.Worksheets.Item("CGC-W-05222025")
Add-TeeFairwayColumn $ws2 $false
$ws2.Range("F13").Select()

$ws3 = $wb.Worksheets.Item("CGC-W-05232025")
Add-TeeFairwayColumn $ws3 $true
$ws3.Range("F2:F10").Select()

$ws4 = $wb.Worksheets.Item("CGC-W-05302025")
Add-TeeFairwayColumn $ws4 $true
$ws4.Range("F2:F10").Select()

$ws5 = $wb.Worksheets.Item("CGC-W-05312025")
Add-TeeFairwayColumn $ws5 $true
$ws5.Range("F2:F10").Select()

# ---------------------------------------------------------------------------
# 3) New sheet for the round played 2025-06-02
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws6.Name = "CGC-W-06022025"

$ws6.Cells.Item(1, 1).Value = "Hole"
$ws6.Cells.Item(1, 2).Value = "Score"
$ws6.Cells.Item(1, 3).Value = "Tee Fairway"
$ws6.Cells.Item(1, 4).Value = "Fairway Hits"
$ws6.Cells.Item(1, 5).Value = "Chips"
$ws6.Cells.Item(1, 6).Value = "Putts"

$rows = @(
    @(1, 9, "Yes", 0, 1, 2),
    @(2, 7, "Yes", 0, 0, 2),
    @(3, 4, "No",  0, 1, 2),
    @(4, 8, "Yes", 0, 4, 2),
    @(5, 6, "Yes", 1, 2, 2),
    @(6, 7, "No",  1, 1, 2),
    @(7, 4, "No",  0, 1, 2),
    @(8, 7, "No",  1, 1, 2),
    @(9, 8, "No",  2, 2, 2)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws6.Cells.Item($r, 1).Value = $row[0]
    $ws6.Cells.Item($r, 2).Value = $row[1]
    $ws6.Cells.Item($r, 3).Value = $row[2]
    $ws6.Cells.Item($r, 4).Value = $row[3]
    $ws6.Cells.Item($r, 5).Value = $row[4]
    $ws6.Cells.Item($r, 6).Value = $row[5]
}

$ws6.Range("D11").Select()
$ws6.Activate()
